$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (pushes Real Madrid / Chelsea FC rows down)
$ws.Rows.Item(3).Insert()

# Row 2: Bayern Munich vs Club Brugge KV - final result 4:0, confidence 75, result column checked
$ws.Range("A2").Value = "Bayern Munich " + [char]0x2713 + " - Club Brugge KV: 4:0"
$ws.Range("C2").Value = 75
$ws.Range("G2").Value = [char]0x2713

# Row 3 (new): Galatasaray vs FK Bodø/Glimt - final result 3:1
# Copy the (empty) Oddspedia_Confidence cell down from row 2 so the new
# row keeps a real, empty placeholder cell in column E, matching the
# other rows' untouched E cells.
$ws.Range("E2").Copy($ws.Range("E3"))
$ws.Range("A3").Value = "Galatasaray " + [char]0x2713 + " - FK Bod" + [char]0x00F8 + "/Glimt: 3:1"
$ws.Range("B3").Value = "Galatasaray"
$ws.Range("C3").Value = 74
$ws.Range("D3").Value = 80
$ws.Range("F3").Value = 1.6
$ws.Range("G3").Value = [char]0x2713

# Row 4: Real Madrid vs Juventus FC - final result 1:0
$ws.Range("A4").Value = "Real Madrid " + [char]0x2713 + " - Juventus FC: 1:0"
$ws.Range("B4").Value = "Real Madrid"
$ws.Range("C4").Value = 72
$ws.Range("D4").Value = 94
$ws.Range("F4").Value = 1.57
$ws.Range("G4").Value = [char]0x2713

# Row 5: Chelsea FC vs Ajax Amsterdam - final result 5:1
$ws.Range("A5").Value = "Chelsea FC " + [char]0x2713 + " - Ajax Amsterdam: 5:1"
$ws.Range("B5").Value = "Chelsea FC"
$ws.Range("C5").Value = 70
$ws.Range("D5").Value = 94
$ws.Range("F5").Value = 1.3
$ws.Range("G5").Value = [char]0x2713
